$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.745.53"
$ws.Range("E2").Value = "  +1.35%  "

$ws.Range("D3").Value = "3.583.82"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'588.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.69%  "

$ws.Range("D6").Value = "'187.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.39%  "

$ws.Range("D7").Value = "3.572.46"
$ws.Range("E7").Value = "  +0.63%  "

$ws.Range("D8").Value = "'0.623"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").Value = "'0.202"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.98%  "

$ws.Range("D11").Value = "'0.652"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.24%  "

$ws.Range("D12").Value = "'54.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.25%  "

$ws.Range("E13").Value = "  +4.66%  "

$ws.Range("D14").Value = "'9.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("D15").Value = "4.155.80"
$ws.Range("E15").Value = "  +0.85%  "

$ws.Range("D16").Value = "'19.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.33%  "

$ws.Range("D17").Value = "70.749.81"
$ws.Range("E17").Value = "  +1.56%  "

$ws.Range("D18").Value = "3.570.40"
$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("D19").Value = "'12.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("D20").Value = "'563.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +14.94%  "

$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'1.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("D23").Value = "'17.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.20%  "

$ws.Range("D24").Value = "'4.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.08%  "

$ws.Range("D25").Value = "'4.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.88%  "

$ws.Range("D26").Value = "'96.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("D27").Value = "'11.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.33%  "

$ws.Range("D28").Value = "'2.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.10%  "

$ws.Range("D29").Value = "'9.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("D30").Value = "'32.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.70%  "

$ws.Range("D31").Value = "'7.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.88%  "

$ws.Range("D32").Value = "'12.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.81%  "

$ws.Range("D33").Value = "'65.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.67%  "

$ws.Range("D34").Value = "'0.115"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.03%  "

$ws.Range("D35").Value = "'566.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("D36").Value = "'3.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.38%  "

$ws.Range("D37").Value = "'0.415"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.67%  "

$ws.Range("D38").Value = "'38.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.82%  "

$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("D40").Value = "0.0₃0778"
$ws.Range("E40").Value = "  -1.00%  "

$ws.Range("E41").Value = "  +1.49%  "

$ws.Range("D42").Value = "3.352.42"
$ws.Range("E42").Value = "  +4.66%  "

$ws.Range("D43").Value = "'3.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.34%  "

$ws.Range("D44").Value = "'3.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.89%  "

$ws.Range("D45").Value = "'3.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.94%  "

$ws.Range("D46").Value = "'2.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.87%  "

$ws.Range("D47").Value = "'0.0446"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.65%  "

$ws.Range("D48").Value = "'9.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.69%  "

$ws.Range("E49").Value = "  +1.60%  "

$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("D51").Value = "'1.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +20.74%  "
